# Updated cryptos list on Mon Feb 12 15:42:09 UTC 2024 with GitHub Actions
#
# This script updates the Price (column D) and Volume(1h) (column E) values
# for the crypto rows on the active sheet, and fixes the FirstDigitalUSD /
# Filecoin row ordering (rows 34-35), to match the freshly scraped data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while forcing Excel to keep it as
# text (it would otherwise silently reinterpret things like "13.07" or
# "2.93" as numbers and lose trailing zeros / introduce float rounding).
# Applying a temporary "@" (Text) number format forces the literal text to
# be stored, then restoring the cell's original Style keeps the workbook's
# style table / formatting untouched.
function Set-TextValue {
    param($Cell, [string]$Value)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "49.379.90"
Set-TextValue $ws.Range("E2") "  +2.60%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.545.71"
Set-TextValue $ws.Range("E3") "  +1.50%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.19%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "321.45"
Set-TextValue $ws.Range("E5") "  +0.27%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "108.21"
Set-TextValue $ws.Range("E6") "  -0.64%  "

# Row 7 - XRP
Set-TextValue $ws.Range("E7") "  -0.76%  "

# Row 8 - USDC
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.13%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.551"
Set-TextValue $ws.Range("E9") "  +1.35%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "40.01"
Set-TextValue $ws.Range("E10") "  +0.37%  "

# Row 11 - Chainlink
Set-TextValue $ws.Range("D11") "20.33"
Set-TextValue $ws.Range("E11") "  +0.83%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("D12") "0.0811"
Set-TextValue $ws.Range("E12") "  -0.76%  "

# Row 13 - TRON
Set-TextValue $ws.Range("E13") "  +0.96%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.21"
Set-TextValue $ws.Range("E14") "  +0.22%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.937.46"
Set-TextValue $ws.Range("E15") "  +1.24%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.549.61"
Set-TextValue $ws.Range("E16") "  +1.97%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.853"
Set-TextValue $ws.Range("E17") "  +0.86%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "49.133.19"
Set-TextValue $ws.Range("E18") "  +2.44%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D19") "13.07"
Set-TextValue $ws.Range("E19") "  -0.89%  "

# Row 20 - ImmutableX
Set-TextValue $ws.Range("D20") "2.93"
Set-TextValue $ws.Range("E20") "  +7.55%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.68"
Set-TextValue $ws.Range("E21") "  +1.23%  "

# Row 22 - ShibaInu
Set-TextValue $ws.Range("D22") "0.0₃0942"
Set-TextValue $ws.Range("E22") "  +0.08%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "283.12"
Set-TextValue $ws.Range("E23") "  +3.36%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "71.49"
Set-TextValue $ws.Range("E24") "  -0.94%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("E25") "  -1.48%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("D26") "26.23"
Set-TextValue $ws.Range("E26") "  +1.35%  "

# Row 27 - Dai
Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  -0.15%  "

# Row 28 - Toncoin
Set-TextValue $ws.Range("D28") "2.23"
Set-TextValue $ws.Range("E28") "  -7.32%  "

# Row 29 - Kaspa
Set-TextValue $ws.Range("D29") "0.145"
Set-TextValue $ws.Range("E29") "  +2.94%  "

# Row 30 - Cosmos
Set-TextValue $ws.Range("D30") "9.78"
Set-TextValue $ws.Range("E30") "  -2.97%  "

# Row 31 - InjectiveProtocol
Set-TextValue $ws.Range("D31") "35.25"
Set-TextValue $ws.Range("E31") "  -0.48%  "

# Row 32 - OKB
Set-TextValue $ws.Range("D32") "49.63"
Set-TextValue $ws.Range("E32") "  +0.34%  "

# Row 33 - Celestia
Set-TextValue $ws.Range("D33") "19.60"
Set-TextValue $ws.Range("E33") "  +1.38%  "

# Row 34 - now Filecoin (was FirstDigitalUSD)
Set-TextValue $ws.Range("B34") "Filecoin"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "5.35"
Set-TextValue $ws.Range("E34") "  +0.17%  "

# Row 35 - now FirstDigitalUSD (was Filecoin)
Set-TextValue $ws.Range("B35") "FirstDigitalUSD"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  -0.36%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("D36") "0.0781"
Set-TextValue $ws.Range("E36") "  -0.19%  "

# Row 37 - ARBITRUM
Set-TextValue $ws.Range("D37") "2.01"
Set-TextValue $ws.Range("E37") "  +2.47%  "

# Row 38 - RenderToken
Set-TextValue $ws.Range("D38") "4.65"
Set-TextValue $ws.Range("E38") "  +0.58%  "

# Row 39 - LidoDAOToken
Set-TextValue $ws.Range("D39") "2.96"
Set-TextValue $ws.Range("E39") "  -0.19%  "

# Row 40 - Stellar
Set-TextValue $ws.Range("E40") "  -0.32%  "

# Row 41 - WEMIXToken
Set-TextValue $ws.Range("D41") "2.22"
Set-TextValue $ws.Range("E41") "  +0.17%  "

# Row 42 - EnergySwap
Set-TextValue $ws.Range("D42") "22.08"
Set-TextValue $ws.Range("E42") "  +1.74%  "

# Row 43 - Monero
Set-TextValue $ws.Range("D43") "119.93"
Set-TextValue $ws.Range("E43") "  -2.05%  "

# Row 44 - VeChain
Set-TextValue $ws.Range("D44") "0.0307"
Set-TextValue $ws.Range("E44") "  +0.78%  "

# Row 45 - NEARProtocol
Set-TextValue $ws.Range("D45") "3.26"
Set-TextValue $ws.Range("E45") "  +4.22%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.013.19"
Set-TextValue $ws.Range("E46") "  -0.33%  "

# Row 47 - Stacks
Set-TextValue $ws.Range("D47") "2.00"
Set-TextValue $ws.Range("E47") "  +7.75%  "

# Row 48 - ApeXProtocol
Set-TextValue $ws.Range("D48") "2.11"
Set-TextValue $ws.Range("E48") "  +6.35%  "

# Row 49 - FraxShare
Set-TextValue $ws.Range("D49") "9.03"
Set-TextValue $ws.Range("E49") "  +0.03%  "

# Row 50 - THORChain
Set-TextValue $ws.Range("D50") "5.27"
Set-TextValue $ws.Range("E50") "  +1.71%  "

# Row 51 - BitcoinSV
Set-TextValue $ws.Range("D51") "80.89"
Set-TextValue $ws.Range("E51") "  +1.76%  "
